$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "59×25=1475" "98×99=9702"
Replace-Text "84×51=4284" "69×63=4347"
Replace-Text "66×84=5544" "38×70=2660"
Replace-Text "90×77=6930" "95×52=4940"
Replace-Text "33×50=1650" "28×64=1792"
Replace-Text "71×25=1775" "27×47=1269"
Replace-Text "72×77=5544" "24×73=1752"
Replace-Text "44×87=3828" "29×41=1189"
Replace-Text "87×34=2958" "48×12=576"
Replace-Text "36×98=3528" "81×62=5022"
Replace-Text "53×88=4664" "28×68=1904"
Replace-Text "77×67=5159" "45×14=630"
Replace-Text "61×54=3294" "14×97=1358"
Replace-Text "30×54=1620" "22×80=1760"
Replace-Text "72×18=1296" "45×17=765"
Replace-Text "83×44=3652" "41×34=1394"
Replace-Text "97×65=6305" "96×84=8064"
Replace-Text "95×30=2850" "98×34=3332"
Replace-Text "86×62=5332" "57×95=5415"
Replace-Text "14×14=196" "89×30=2670"
Replace-Text "99×85=8415" "75×85=6375"
Replace-Text "69×68=4692" "76×86=6536"
Replace-Text "16×16=256" "97×22=2134"
Replace-Text "91×93=8463" "71×88=6248"
Replace-Text "81×63=5103" "28×29=812"
